# 自动更新Excel文件 - 2026-02-10 23:32:08
# Recompute the "剩余" (remaining days) column (E) for every data row based on the
# current reference date (2026-02-11). Remaining = (开始时间(F) + 总天(D)) - today.
# If a cycle has already reached/passed its end date, it is renewed: the start
# date (F) rolls forward to the old end date and the remaining days resets to
# the full duration (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = [datetime]::ParseExact("20260211", "yyyyMMdd", $null)
$todayOA = $today.ToOADate()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)   # D: 总天 (total days)
    $eCell = $ws.Cells.Item($r, 5)   # E: 剩余 (remaining days)
    $fCell = $ws.Cells.Item($r, 6)   # F: 开始时间 (start date, yyyyMMdd)

    $dVal = $dCell.Value2
    $fVal = $fCell.Value2

    if (-not $dVal) { continue }
    if (-not $fVal) { continue }

    $fStr = "{0:0}" -f $fVal

    $startDate = $null
    try {
        $startDate = [datetime]::ParseExact($fStr, "yyyyMMdd", $null)
    } catch {
        $startDate = $null
    }

    if (-not $startDate) { continue }

    $endDate = $startDate.AddDays($dVal)
    $remaining = [math]::Round($endDate.ToOADate() - $todayOA)

    if ($remaining -le 0) {
        # Cycle finished - renew it starting from the old end date.
        $startDate = $endDate
        $endDate = $startDate.AddDays($dVal)
        $remaining = [math]::Round($endDate.ToOADate() - $todayOA)

        $fCell.Value2 = [double]($startDate.ToString("yyyyMMdd"))
    }

    $eCell.Value2 = $remaining
}
